# A new "Others" entry was logged for September on 2024-09-03 19:13:40
# ("login internet personal share"). The tracker sheet keeps each month's
# entries stacked at the top of their section, so this pushes the existing
# rows 29-42 (of the "2024" sheet) down by one row to row 43.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row above row 29, shifting rows 29:42 down to 30:43.
$ws.Rows(29).Insert()

# Populate the newly inserted row with the new September entry.
$ws.Range("R29").Value = "login internet personal share"
$ws.Range("S29").Value = "2024-09-03 19:13:40"
